# docs(taskplan): update task management for sigin up, sign in
#
# Fill in the start/end dates for the two "sign up / sign in" related tasks
# on the "Task List" sheet (rows 6 and 7, columns D = start date, E = end date).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task List")

# Row 6: "Đăng ký/Đăng nhập/Đăng xuất" (Sign up/Sign in/Sign out)
$ws.Range("D6").Value = 45827
$ws.Range("E6").Value = 45831

# Row 7: "Route bảo vệ, test case Auth" (Protected routes, Auth test cases)
$ws.Range("D7").Value = 45827
$ws.Range("E7").Value = 45831
